# Updates the bulk user credential sheet with newly generated iAuthor TC credentials.
# For each data row (2-15) the Client Id, User Name, Exam Password, First Name and
# Last Name columns get fresh randomised values, and the Candidate ID (col B) is bumped
# to the new numbering scheme (231102226-231102239).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 'ktHIK344'
$ws.Range("B2").Value = 231102239
$ws.Range("C2").Value = 'mbjcekk23'
$ws.Range("D2").Value = 'T%rA#45d'
$ws.Range("F2").Value = 'JGXPVFKr'
$ws.Range("G2").Value = 'DjTJ'

# Row 3
$ws.Range("A3").Value = 'SQHYx858'
$ws.Range("B3").Value = 231102238
$ws.Range("C3").Value = 'welkjhn58'
$ws.Range("D3").Value = 'U5h%D$k2'
$ws.Range("F3").Value = 'EmNmkpZU'
$ws.Range("G3").Value = 'yZSB'

# Row 4
$ws.Range("A4").Value = 'YHNhv302'
$ws.Range("B4").Value = 231102237
$ws.Range("C4").Value = 'oezfisi13'
$ws.Range("D4").Value = 'WQ62#u&j'
$ws.Range("F4").Value = 'GcMoaubA'
$ws.Range("G4").Value = 'jvzK'

# Row 5
$ws.Range("A5").Value = 'bELQv496'
$ws.Range("B5").Value = 231102236
$ws.Range("C5").Value = 'wxfvjml28'
$ws.Range("D5").Value = 'w&8G4yR%'
$ws.Range("F5").Value = 'ePQinYrz'
$ws.Range("G5").Value = 'Wnvo'

# Row 6
$ws.Range("A6").Value = 'gqUYb661'
$ws.Range("B6").Value = 231102235
$ws.Range("C6").Value = 'lxbqbul57'
$ws.Range("D6").Value = 'A5#$Wv6h'
$ws.Range("F6").Value = 'DlXxbdiw'
$ws.Range("G6").Value = 'ocvZ'

# Row 7
$ws.Range("A7").Value = 'giMXD898'
$ws.Range("B7").Value = 231102234
$ws.Range("C7").Value = 'kbzjzrr65'
$ws.Range("D7").Value = 'x&5eH!3B'
$ws.Range("F7").Value = 'ChUVPJGA'
$ws.Range("G7").Value = 'LXmJ'

# Row 8
$ws.Range("A8").Value = 'LxfXa374'
$ws.Range("B8").Value = 231102233
$ws.Range("C8").Value = 'jcguvpl47'
$ws.Range("D8").Value = 'ZP%5!q4s'
$ws.Range("F8").Value = 'oDRpqPKu'
$ws.Range("G8").Value = 'JGiX'

# Row 9
$ws.Range("A9").Value = 'cniQj945'
$ws.Range("B9").Value = 231102232
$ws.Range("C9").Value = 'yqbwwyl98'
$ws.Range("D9").Value = 'vJq8#K2%'
$ws.Range("F9").Value = 'DRdaRSrl'
$ws.Range("G9").Value = 'MBJo'

# Row 10
$ws.Range("A10").Value = 'ZSLCT519'
$ws.Range("B10").Value = 231102231
$ws.Range("C10").Value = 'muojxoy68'
$ws.Range("D10").Value = 'Ej$7&pA5'
$ws.Range("F10").Value = 'LHGsxFEd'
$ws.Range("G10").Value = 'TzOx'

# Row 11
$ws.Range("A11").Value = 'zMgTf231'
$ws.Range("B11").Value = 231102230
$ws.Range("C11").Value = 'jhfbvyd52'
$ws.Range("D11").Value = 'dV%&bA49'
$ws.Range("F11").Value = 'PcxCjfLZ'
$ws.Range("G11").Value = 'vgbS'

# Row 12
$ws.Range("A12").Value = 'ZTItV653'
$ws.Range("B12").Value = 231102229
$ws.Range("C12").Value = 'nwstiba64'
$ws.Range("D12").Value = 'SY5&m9!a'
$ws.Range("F12").Value = 'NTVmTljh'
$ws.Range("G12").Value = 'vlei'

# Row 13
$ws.Range("A13").Value = 'XnJXG624'
$ws.Range("B13").Value = 231102228
$ws.Range("C13").Value = 'lkgifoz18'
$ws.Range("D13").Value = 'K3$#2cXd'
$ws.Range("F13").Value = 'DuYniwes'
$ws.Range("G13").Value = 'LYty'

# Row 14
$ws.Range("A14").Value = 'WrHkh461'
$ws.Range("B14").Value = 231102227
$ws.Range("C14").Value = 'tknpvxn89'
$ws.Range("D14").Value = 'VZ$2g9s%'
$ws.Range("F14").Value = 'rCCvALWJ'
$ws.Range("G14").Value = 'XHTL'

# Row 15
$ws.Range("A15").Value = 'RunVP388'
$ws.Range("B15").Value = 231102226
$ws.Range("C15").Value = 'bwlfbth62'
$ws.Range("D15").Value = 'PWf&6s7%'
$ws.Range("F15").Value = 'YboUEKNX'
$ws.Range("G15").Value = 'yKvy'

